$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "322. Coin Change"
$ws.Range("B3").Value = 'recursive, dp, each node returns "least number of coins needed to build its value". At each node take lleast out of all child and + 1 '

$ws.Range("A4").Value = "238. Product of Array Except Self"
$ws.Range("B4").Value = "make container, populate from left to right and right to left O(2n). Value at each index of container = container[i-1] * nums[i-1]"

$null = $ws.Range("B4").Select()
